# Update the "timestamp" column (column O) for rows 2 through 73
# from "2022-07-20 07:01:57" to "2022-07-20 20:58:53".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 73; $row++) {
    $cell = $ws.Cells.Item($row, 15)  # Column O is the 15th column
    if ($cell.Value2 -eq "2022-07-20 07:01:57") {
        $cell.Value2 = "2022-07-20 20:58:53"
    }
}
